$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 184.0626906666667
$ws.Range("H2").Value = 552.188072
$ws.Range("I2").Value = 0.6510505751503485
$ws.Range("J2").Value = 0.6510505751503486
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8749903333333333
$ws.Range("N2").Value = 2.624971
$ws.Range("O2").Value = 0.2670516933349977
$ws.Range("P2").Value = 0.2670516933349977
$ws.Range("Q2").Value = 161.0530750606569
$ws.Range("R2").Value = 1449.477675545912
$ws.Range("S2").Value = 0.1738641585406248
$ws.Range("T2").Value = 0.1738641585406248

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 184.0626906666667
$ws.Range("H3").Value = 552.188072
$ws.Range("I3").Value = 0.6510505751503485
$ws.Range("J3").Value = 0.6510505751503486
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.8147036666666666
$ws.Range("N3").Value = 2.444111
$ws.Range("O3").Value = 0.2486518827250642
$ws.Range("P3").Value = 0.2486518827250642
$ws.Range("Q3").Value = 149.9565489826658
$ws.Range("R3").Value = 1349.608940843992
$ws.Range("S3").Value = 0.1618849512603701
$ws.Range("T3").Value = 0.1618849512603701

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 184.0626906666667
$ws.Range("H4").Value = 552.188072
$ws.Range("I4").Value = 0.6510505751503485
$ws.Range("J4").Value = 0.6510505751503486
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.586789
$ws.Range("N4").Value = 4.760367
$ws.Range("O4").Value = 0.484296423939938
$ws.Range("P4").Value = 0.484296423939938
$ws.Range("Q4").Value = 292.0686528602693
$ws.Range("R4").Value = 2628.617875742424
$ws.Range("S4").Value = 0.3153014653493536
$ws.Range("T4").Value = 0.3153014653493537

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 57.4434
$ws.Range("H5").Value = 172.3302
$ws.Range("I5").Value = 0.2031838091312023
$ws.Range("J5").Value = 0.2031838091312023
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.8749903333333333
$ws.Range("N5").Value = 2.624971
$ws.Range("O5").Value = 0.2670516933349977
$ws.Range("P5").Value = 0.2670516933349977
$ws.Range("Q5").Value = 50.26241971379999
$ws.Range("R5").Value = 452.3617774241999
$ws.Range("S5").Value = 0.05426058028674255
$ws.Range("T5").Value = 0.05426058028674255

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 57.4434
$ws.Range("H6").Value = 172.3302
$ws.Range("I6").Value = 0.2031838091312023
$ws.Range("J6").Value = 0.2031838091312023
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.8147036666666666
$ws.Range("N6").Value = 2.444111
$ws.Range("O6").Value = 0.2486518827250642
$ws.Range("P6").Value = 0.2486518827250642
$ws.Range("Q6").Value = 46.79934860579999
$ws.Range("R6").Value = 421.1941374522
$ws.Range("S6").Value = 0.05052203667972355
$ws.Range("T6").Value = 0.05052203667972355

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 57.4434
$ws.Range("H7").Value = 172.3302
$ws.Range("I7").Value = 0.2031838091312023
$ws.Range("J7").Value = 0.2031838091312023
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.586789
$ws.Range("N7").Value = 4.760367
$ws.Range("O7").Value = 0.484296423939938
$ws.Range("P7").Value = 0.484296423939938
$ws.Range("Q7").Value = 91.15055524259998
$ws.Range("R7").Value = 820.3549971833999
$ws.Range("S7").Value = 0.0984011921647362
$ws.Range("T7").Value = 0.0984011921647362

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 41.21033366666666
$ws.Range("H8").Value = 123.631001
$ws.Range("I8").Value = 0.1457656157184491
$ws.Range("J8").Value = 0.1457656157184491
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.8749903333333333
$ws.Range("N8").Value = 2.624971
$ws.Range("O8").Value = 0.2670516933349977
$ws.Range("P8").Value = 0.2670516933349977
$ws.Range("Q8").Value = 36.05864359177455
$ws.Range("R8").Value = 324.527792325971
$ws.Range("S8").Value = 0.0389269545076304
$ws.Range("T8").Value = 0.0389269545076304

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 41.21033366666666
$ws.Range("H9").Value = 123.631001
$ws.Range("I9").Value = 0.1457656157184491
$ws.Range("J9").Value = 0.1457656157184491
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.8147036666666666
$ws.Range("N9").Value = 2.444111
$ws.Range("O9").Value = 0.2486518827250642
$ws.Range("P9").Value = 0.2486518827250642
$ws.Range("Q9").Value = 33.57420994279011
$ws.Range("R9").Value = 302.167889485111
$ws.Range("S9").Value = 0.0362448947849706
$ws.Range("T9").Value = 0.03624489478497059

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 41.21033366666666
$ws.Range("H10").Value = 123.631001
$ws.Range("I10").Value = 0.1457656157184491
$ws.Range("J10").Value = 0.1457656157184491
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.586789
$ws.Range("N10").Value = 4.760367
$ws.Range("O10").Value = 0.484296423939938
$ws.Range("P10").Value = 0.484296423939938
$ws.Range("Q10").Value = 65.39210414859632
$ws.Range("R10").Value = 588.528937337367
$ws.Range("S10").Value = 0.07059376642584812
$ws.Range("T10").Value = 0.07059376642584812
